# Regenerate orders with updated distance/size codes.
# Distances: D51 -> D55, D64 -> D69, D80 -> D86
# Sizes:     S30 -> S31
# These codes appear inside many shared strings (Condition, Filename_Left,
# Filename_Right, Distance, Size columns), so a sheet-wide text Replace is
# the most direct way to reproduce the rename across every occurrence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.Cells

# Order matters only in that these four tokens are mutually exclusive
# (no string contains more than one), so a straightforward sequential
# replace is safe - none of the replacements can collide with another.
$cells.Replace("D51", "D55")
$cells.Replace("D64", "D69")
$cells.Replace("D80", "D86")
$cells.Replace("S30", "S31")
